$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update song name (B2): "Heya" -> "Sone Lagde"
$ws.Range("B2").Value = "Sone Lagde"

# Update song URL (C2 and C3) to the new pagalworld.com.bz link
$ws.Range("C2").Value = "https://www.pagalworld.com.bz/dl/miss-pooja-sohne-lagde-mp3-song-download/64"
$ws.Range("C3").Value = "https://www.pagalworld.com.bz/dl/miss-pooja-sohne-lagde-mp3-song-download/64"

# Remove the hyperlinks that pointed at the old URL
$ws.Hyperlinks.Delete()

# Move the view/selection to reflect the saved state (top-left visible cell B1, active cell C3)
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("C3").Select()
